# The trailing empty paragraph (currently <w:p/>, with no runs at all) gains
# an empty run: <w:p><w:r><w:t/></w:r></w:p>. Inserting an empty string into
# that paragraph's range materializes the (empty) run/text without changing
# the visible content or adding a new paragraph.
$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("")
